# Add "Library pin" (column B) values for the three MCP23017 pin-map
# tables on Sheet1, and update the saved cursor/selection position to
# reflect where the author was last working (C61).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The same 16-row Library-pin pattern repeats for each of the three
# MCP23017 IC tables (data rows 36-51, 55-70, 74-89).
$libraryPins = @(8, 9, 10, 11, 12, 13, 14, 15, 0, 1, 2, 3, 4, 5, 6, 7)

$blockStartRows = @(36, 55, 74)

foreach ($startRow in $blockStartRows) {
    for ($i = 0; $i -lt $libraryPins.Length; $i++) {
        $row = $startRow + $i
        $ws.Cells.Item($row, 2).Value = $libraryPins[$i]
    }
}

# Update the view: scroll so row 61 is the top-left visible row, and
# select C61 as the active cell.
$ws.Application.ActiveWindow.ScrollRow = 61
$ws.Range("C61").Select()
